$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AutoCompleteSampleSheet")
$ws.Activate()
$ws.Range("A2").Value = "JavaScript"
$ws.Range("A10").Select()
